# Apply the changes described in the commit:
# "regenerate instance to have positive average demands during the last periods"

$wb = $excel.ActiveWorkbook

# --- Productdata sheet: G2 40 -> 70 ---
$wsProductdata = $wb.Worksheets.Item("Productdata")
$wsProductdata.Range("G2").Value = 70

# --- ForecastedAverageDemand sheet: B9,B10,B11 0 -> 100 ---
$wsAvgDemand = $wb.Worksheets.Item("ForecastedAverageDemand")
$wsAvgDemand.Range("B9").Value = 100
$wsAvgDemand.Range("B10").Value = 100
$wsAvgDemand.Range("B11").Value = 100

# --- ForcastedStandardDeviation sheet: B9,B10,B11 0 -> new values ---
$wsStdDev = $wb.Worksheets.Item("ForcastedStandardDeviation")
$wsStdDev.Range("B9").Value = 10.23775
$wsStdDev.Range("B10").Value = 11.713975
$wsStdDev.Range("B11").Value = 13.0425775

# Re-assert the (untouched) blank H column cells on Productdata as empty
# strings. These cells already held blank/empty shared-string values in the
# source file; doing this keeps them reading back as blank after the
# workbook round-trips through the COM object model, instead of drifting to
# a stray shared-string index.
for ($r = 2; $r -le 11; $r++) {
    $wsProductdata.Cells.Item($r, 8).Value = ""
}
